$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 7 with the new faculty contact entry (Sini Mam / Chemistry)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 45126
$ws.Range("C7").Value = "Sini Mam"
$ws.Range("D7").Value = "Chemistry"
$ws.Range("E7").Value = 8240897581
$ws.Range("F7").Value = "Vikramjit Chakraborty (Gullu son of Abhishek Chakraborty friend of Anirban)"

# Update the active selection to H7, matching the saved view state
$ws.Range("H7").Select()
